$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 37.17328633333333
$ws.Range("H2").Value = 111.519859
$ws.Range("I2").Value = 0.005170079968594893
$ws.Range("J2").Value = 0.005188590814393131
$ws.Range("M2").Value = 51.02156433333334
$ws.Range("N2").Value = 153.064693
$ws.Range("O2").Value = 0.4760900215891154
$ws.Range("P2").Value = 0.4807539937572116
$ws.Range("Q2").Value = 1896.639220137587
$ws.Range("R2").Value = 17069.75298123829
$ws.Range("S2").Value = 0.002461423483865795
$ws.Range("T2").Value = 0.002494435755991481
$ws.Range("G3").Value = 37.17328633333333
$ws.Range("H3").Value = 111.519859
$ws.Range("I3").Value = 0.005170079968594893
$ws.Range("J3").Value = 0.005188590814393131
$ws.Range("O3").Value = 0.008900168787493621
$ws.Range("P3").Value = 0.008987358473548528
$ws.Range("Q3").Value = 35.45633897526477
$ws.Range("R3").Value = 319.1070507773829
$ws.Range("S3").Value = 0.00004601458436533426
$ws.Range("T3").Value = 0.00004663172562151216
$ws.Range("G4").Value = 37.17328633333333
$ws.Range("H4").Value = 111.519859
$ws.Range("I4").Value = 0.005170079968594893
$ws.Range("J4").Value = 0.005188590814393131
$ws.Range("M4").Value = 28.11170133333333
$ws.Range("N4").Value = 84.335104
$ws.Range("O4").Value = 0.2623145854026591
$ws.Range("P4").Value = 0.2648843261452188
$ws.Range("Q4").Value = 1045.004322981148
$ws.Range("R4").Value = 9405.038906830336
$ws.Range("S4").Value = 0.001356187383460562
$ws.Range("T4").Value = 0.001374376381513797
$ws.Range("G5").Value = 37.17328633333333
$ws.Range("H5").Value = 111.519859
$ws.Range("I5").Value = 0.005170079968594893
$ws.Range("J5").Value = 0.005188590814393131
$ws.Range("M5").Value = 3.119026
$ws.Range("N5").Value = 6.238052
$ws.Range("O5").Value = 0.02910410872500189
$ws.Range("P5").Value = 0.01959281630196169
$ws.Range("Q5").Value = 115.9444465791113
$ws.Range("R5").Value = 695.6666794746679
$ws.Range("S5").Value = 0.0001504705695229401
$ws.Range("T5").Value = 0.0001016591066924504
$ws.Range("G6").Value = 37.17328633333333
$ws.Range("H6").Value = 111.519859
$ws.Range("I6").Value = 0.005170079968594893
$ws.Range("J6").Value = 0.005188590814393131
$ws.Range("M6").Value = 23.96178866666667
$ws.Range("N6").Value = 71.885366
$ws.Range("O6").Value = 0.22359111549573
$ws.Range("P6").Value = 0.2257815053220593
$ws.Range("Q6").Value = 890.7384311648215
$ws.Range("R6").Value = 8016.645880483395
$ws.Range("S6").Value = 0.001155983947380261
$ws.Range("T6").Value = 0.001171487844573891
$ws.Range("I7").Value = 0.006280726092526873
$ws.Range("J7").Value = 0.006303213472394487
$ws.Range("M7").Value = 51.02156433333334
$ws.Range("N7").Value = 153.064693
$ws.Range("O7").Value = 0.4760900215891154
$ws.Range("P7").Value = 0.4807539937572116
$ws.Range("Q7").Value = 2304.078759011042
$ws.Range("R7").Value = 20736.70883109937
$ws.Range("S7").Value = 0.002990191020986439
$ws.Range("T7").Value = 0.003030295050357912
$ws.Range("I8").Value = 0.006280726092526873
$ws.Range("J8").Value = 0.006303213472394487
$ws.Range("O8").Value = 0.008900168787493621
$ws.Range("P8").Value = 0.008987358473548528
$ws.Range("S8").Value = 0.00005589952233150445
$ws.Range("T8").Value = 0.00005664923901170984
$ws.Range("I9").Value = 0.006280726092526873
$ws.Range("J9").Value = 0.006303213472394487
$ws.Range("M9").Value = 28.11170133333333
$ws.Range("N9").Value = 84.335104
$ws.Range("O9").Value = 0.2623145854026591
$ws.Range("P9").Value = 0.2648843261452188
$ws.Range("Q9").Value = 1269.49408094646
$ws.Range("R9").Value = 11425.44672851815
$ws.Range("S9").Value = 0.001647526060988849
$ws.Range("T9").Value = 0.001669622453184679
$ws.Range("I10").Value = 0.006280726092526873
$ws.Range("J10").Value = 0.006303213472394487
$ws.Range("M10").Value = 3.119026
$ws.Range("N10").Value = 6.238052
$ws.Range("O10").Value = 0.02910410872500189
$ws.Range("P10").Value = 0.01959281630196169
$ws.Range("Q10").Value = 140.8518466515953
$ws.Range("R10").Value = 845.111079909572
$ws.Range("S10").Value = 0.0001827949350688584
$ws.Range("T10").Value = 0.0001234977036766753
$ws.Range("I11").Value = 0.006280726092526873
$ws.Range("J11").Value = 0.006303213472394487
$ws.Range("M11").Value = 23.96178866666667
$ws.Range("N11").Value = 71.885366
$ws.Range("O11").Value = 0.22359111549573
$ws.Range("P11").Value = 0.2257815053220593
$ws.Range("Q11").Value = 1082.08850544217
$ws.Range("R11").Value = 9738.796548979528
$ws.Range("S11").Value = 0.001404314553151221
$ws.Range("T11").Value = 0.001423149026163512
$ws.Range("G12").Value = 2375.59786
$ws.Range("H12").Value = 7126.79358
$ws.Range("I12").Value = 0.3303993840977568
$ws.Range("J12").Value = 0.3315823391174117
$ws.Range("M12").Value = 51.02156433333334
$ws.Range("N12").Value = 153.064693
$ws.Range("O12").Value = 0.4760900215891154
$ws.Range("P12").Value = 0.4807539937572116
$ws.Range("Q12").Value = 121206.719044119
$ws.Range("R12").Value = 1090860.471397071
$ws.Range("S12").Value = 0.1572998499081315
$ws.Range("T12").Value = 0.1594095337900538
$ws.Range("G13").Value = 2375.59786
$ws.Range("H13").Value = 7126.79358
$ws.Range("I13").Value = 0.3303993840977568
$ws.Range("J13").Value = 0.3315823391174117
$ws.Range("O13").Value = 0.008900168787493621
$ws.Range("P13").Value = 0.008987358473548528
$ws.Range("Q13").Value = 2265.874537908273
$ws.Range("R13").Value = 20392.87084117446
$ws.Range("S13").Value = 0.002940610285753971
$ws.Range("T13").Value = 0.002980049345145911
$ws.Range("G14").Value = 2375.59786
$ws.Range("H14").Value = 7126.79358
$ws.Range("I14").Value = 0.3303993840977568
$ws.Range("J14").Value = 0.3315823391174117
$ws.Range("M14").Value = 28.11170133333333
$ws.Range("N14").Value = 84.335104
$ws.Range("O14").Value = 0.2623145854026591
$ws.Range("P14").Value = 0.2648843261452188
$ws.Range("Q14").Value = 66782.0975284258
$ws.Range("R14").Value = 601038.8777558323
$ws.Range("S14").Value = 0.08666857745689699
$ws.Range("T14").Value = 0.08783096445877103
$ws.Range("G15").Value = 2375.59786
$ws.Range("H15").Value = 7126.79358
$ws.Range("I15").Value = 0.3303993840977568
$ws.Range("J15").Value = 0.3315823391174117
$ws.Range("M15").Value = 3.119026
$ws.Range("N15").Value = 6.238052
$ws.Range("O15").Value = 0.02910410872500189
$ws.Range("P15").Value = 0.01959281630196169
$ws.Range("Q15").Value = 7409.551490884359
$ws.Range("R15").Value = 44457.30894530615
$ws.Range("S15").Value = 0.009615979597454775
$ws.Range("T15").Value = 0.006496631859302214
$ws.Range("G16").Value = 2375.59786
$ws.Range("H16").Value = 7126.79358
$ws.Range("I16").Value = 0.3303993840977568
$ws.Range("J16").Value = 0.3315823391174117
$ws.Range("M16").Value = 23.96178866666667
$ws.Range("N16").Value = 71.885366
$ws.Range("O16").Value = 0.22359111549573
$ws.Range("P16").Value = 0.2257815053220593
$ws.Range("Q16").Value = 56923.57387830559
$ws.Range("R16").Value = 512312.1649047503
$ws.Range("S16").Value = 0.0738743668495196
$ws.Range("T16").Value = 0.07486515966413877
$ws.Range("G17").Value = 76.954105
$ws.Range("H17").Value = 153.90821
$ws.Range("I17").Value = 0.01070281688829022
$ws.Range("J17").Value = 0.007160758019481436
$ws.Range("M17").Value = 51.02156433333334
$ws.Range("N17").Value = 153.064693
$ws.Range("O17").Value = 0.4760900215891154
$ws.Range("P17").Value = 0.4807539937572116
$ws.Range("Q17").Value = 3926.318818971588
$ws.Range("R17").Value = 23557.91291382953
$ws.Range("S17").Value = 0.005095504323410441
$ws.Range("T17").Value = 0.003442563016194681
$ws.Range("G18").Value = 76.954105
$ws.Range("H18").Value = 153.90821
$ws.Range("I18").Value = 0.01070281688829022
$ws.Range("J18").Value = 0.007160758019481436
$ws.Range("O18").Value = 0.008900168787493621
$ws.Range("P18").Value = 0.008987358473548528
$ws.Range("Q18").Value = 73.39977444962832
$ws.Range("R18").Value = 440.3986466977699
$ws.Range("S18").Value = 0.00009525687680742022
$ws.Range("T18").Value = 0.00006435629926341706
$ws.Range("G19").Value = 76.954105
$ws.Range("H19").Value = 153.90821
$ws.Range("I19").Value = 0.01070281688829022
$ws.Range("J19").Value = 0.007160758019481436
$ws.Range("M19").Value = 28.11170133333333
$ws.Range("N19").Value = 84.335104
$ws.Range("O19").Value = 0.2623145854026591
$ws.Range("P19").Value = 0.2648843261452188
$ws.Range("Q19").Value = 2163.310816133973
$ws.Range("R19").Value = 12979.86489680384
$ws.Range("S19").Value = 0.002807504974692427
$ws.Range("T19").Value = 0.001896772562679312
$ws.Range("G20").Value = 76.954105
$ws.Range("H20").Value = 153.90821
$ws.Range("I20").Value = 0.01070281688829022
$ws.Range("J20").Value = 0.007160758019481436
$ws.Range("M20").Value = 3.119026
$ws.Range("N20").Value = 6.238052
$ws.Range("O20").Value = 0.02910410872500189
$ws.Range("P20").Value = 0.01959281630196169
$ws.Range("Q20").Value = 240.02185430173
$ws.Range("R20").Value = 960.0874172069199
$ws.Range("S20").Value = 0.000311495946380585
$ws.Range("T20").Value = 0.0001402994164584988
$ws.Range("G21").Value = 76.954105
$ws.Range("H21").Value = 153.90821
$ws.Range("I21").Value = 0.01070281688829022
$ws.Range("J21").Value = 0.007160758019481436
$ws.Range("M21").Value = 23.96178866666667
$ws.Range("N21").Value = 71.885366
$ws.Range("O21").Value = 0.22359111549573
$ws.Range("P21").Value = 0.2257815053220593
$ws.Range("Q21").Value = 1843.958001042477
$ws.Range("R21").Value = 11063.74800625486
$ws.Range("S21").Value = 0.002393054766999349
$ws.Range("T21").Value = 0.001616766724885527
$ws.Range("G22").Value = 4655.195393666666
$ws.Range("H22").Value = 13965.586181
$ws.Range("I22").Value = 0.6474469929528313
$ws.Range("J22").Value = 0.6497650985763194
$ws.Range("M22").Value = 51.02156433333334
$ws.Range("N22").Value = 153.064693
$ws.Range("O22").Value = 0.4760900215891154
$ws.Range("P22").Value = 0.4807539937572116
$ws.Range("Q22").Value = 237515.3512622008
$ws.Range("R22").Value = 2137638.161359807
$ws.Range("S22").Value = 0.3082430528527213
$ws.Range("T22").Value = 0.3123771661446138
$ws.Range("G23").Value = 4655.195393666666
$ws.Range("H23").Value = 13965.586181
$ws.Range("I23").Value = 0.6474469929528313
$ws.Range("J23").Value = 0.6497650985763194
$ws.Range("O23").Value = 0.008900168787493621
$ws.Range("P23").Value = 0.008987358473548528
$ws.Range("Q23").Value = 4440.182780555787
$ws.Range("R23").Value = 39961.64502500209
$ws.Range("S23").Value = 0.005762387518235391
$ws.Range("T23").Value = 0.005839671864505978
$ws.Range("G24").Value = 4655.195393666666
$ws.Range("H24").Value = 13965.586181
$ws.Range("I24").Value = 0.6474469929528313
$ws.Range("J24").Value = 0.6497650985763194
$ws.Range("M24").Value = 28.11170133333333
$ws.Range("N24").Value = 84.335104
$ws.Range("O24").Value = 0.2623145854026591
$ws.Range("P24").Value = 0.2648843261452188
$ws.Range("Q24").Value = 130865.4625550664
$ws.Range("R24").Value = 1177789.162995598
$ws.Range("S24").Value = 0.1698347895266203
$ws.Range("T24").Value = 0.1721125902890701
$ws.Range("G25").Value = 4655.195393666666
$ws.Range("H25").Value = 13965.586181
$ws.Range("I25").Value = 0.6474469929528313
$ws.Range("J25").Value = 0.6497650985763194
$ws.Range("M25").Value = 3.119026
$ws.Range("N25").Value = 6.238052
$ws.Range("O25").Value = 0.02910410872500189
$ws.Range("P25").Value = 0.01959281630196169
$ws.Range("Q25").Value = 14519.67546792657
$ws.Range("R25").Value = 87118.0528075594
$ws.Range("S25").Value = 0.01884336767657473
$ws.Range("T25").Value = 0.01273072821583185
$ws.Range("G26").Value = 4655.195393666666
$ws.Range("H26").Value = 13965.586181
$ws.Range("I26").Value = 0.6474469929528313
$ws.Range("J26").Value = 0.6497650985763194
$ws.Range("M26").Value = 23.96178866666667
$ws.Range("N26").Value = 71.885366
$ws.Range("O26").Value = 0.22359111549573
$ws.Range("P26").Value = 0.2257815053220593
$ws.Range("Q26").Value = 111546.8082250808
$ws.Range("R26").Value = 1003921.274025727
$ws.Range("S26").Value = 0.1447633953786796
$ws.Range("T26").Value = 0.1467049420622976
